$d = $word.ActiveDocument

# --- Change 1: split "Spielfeld wird mit 0 initialisiert" and move the _GoBack
#     bookmark into the middle of that run (between "0" and " initialisiert").
$rng = $d.Content
$rng.Find.Execute("Spielfeld wird mit 0")
$insertPos = $rng.End
$bmRange = $d.Range($insertPos, $insertPos)

# --- Change 2: the existing _GoBack bookmark (currently located after the
#     "Snake besitzt eine Exemplarvariable..." paragraph) must be removed from
#     its old spot. Word only ever has a single _GoBack bookmark, so deleting
#     it first (if present) and then adding the new one has the same effect
#     as "moving" it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Change 3: replace the empty paragraph that follows "previousY" with a
#     new list paragraph "Bausteine von Snake" (same list/style as the other
#     top-level bullet items: ListParagraph, ilvl=0, numId=2).
$rng2 = $d.Content
$rng2.Find.Execute("previousY")
$para = $rng2.Paragraphs(1)
$nextPara = $para.Next()
$nextRange = $nextPara.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Bausteine von Snake</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$nextRange.InsertXML($xml)
